# Actualización desde MV -datos-
# Append the next 6 daily dates (09-10-2021 .. 14-10-2021) to the bottom
# of the "Bonos Bancarios en Pesos ($)" daily table, each with the same
# B/C values (322 / 0) as the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("09-10-2021", "10-10-2021", "11-10-2021", "12-10-2021", "13-10-2021", "14-10-2021")

# Last populated row in the current table (column A, walking up from the
# bottom of the sheet - robust even if UsedRange ever grows stale).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$startRow = $lastRow + 1

# Scratch cell used only to get the literal date text into the sheet as a
# plain string (via a text formula + paste-values) so Excel's automatic
# "looks like a date" conversion never kicks in and no cell formatting is
# touched in the process.
$helper = $ws.Range("Z1")

for ($i = 0; $i -lt $dates.Count; $i++) {
    $row = $startRow + $i
    $helper.Formula = '="' + $dates[$i] + '"'
    $helper.Copy()
    $ws.Range("A" + $row).PasteSpecial(-4163)
    $ws.Range("B" + $row).Value = 322
    $ws.Range("C" + $row).Value = 0
}

$helper.Clear()
$excel.CutCopyMode = $false
